$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Custody paid dividend on 12,000 shares while entitlement should be 10,000; the CHF 4,030 net difference equals 2,000 shares × CHF 2.015 net/share."
$ws.Range("H2").Value = 0.86
$ws.Range("I2").Value = "DRAFT_CUSTODIAN_TICKET"
$ws.Range("J2").Value = "Please verify the dividend entitlement quantity at record date. Custody file states NOMINAL_BASIS 10,000 but amounts reflect 12,000 shares. If 2,000 shares were included in error, kindly advise reversal/adjustment; otherwise provide rationale (e.g., late trade or sub-account allocation). Event dates align on both sides (Ex 25-Apr-2025, Pay 29-Apr-2025)."

$ws.Range("G3").Value = "Custody applied 20% tax (net QC 7,220,000 KRW) while NBIM applied ~25% total (22% WHT + ~2.985% local; net QC 6,769,950 KRW), leading to higher custody cash by 342.77 USD. Dates and lending flags do not explain the amount; this is a tax-rate application difference."
$ws.Range("H3").Value = 0.9
$ws.Range("I3").Value = "DRAFT_CUSTODIAN_TICKET"
$ws.Range("J3").Value = "Please ask custodian (CUST/HSBCKR) to confirm the correct Korean dividend withholding and local surtax for this event (Samsung Electronics, ex 31-Mar-2025, pay 20/25-May-2025). Their booking reflects 20% only, while NBIM expects ~25% total (22% WHT + local). Clarify whether local surtax and the additional 2% WHT component will be collected/adjusted, and whether any securities lending portion (2,000 shares) is treated differently. Cross-currency reversal noted but not the driver of the discrepancy."
